$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Importe" column values that were scraped with the wrong
# decimal separator formatting ("836,00" -> "836.00", "8.753,00" ->
# "8753.00"). These cells hold plain text (not real numbers), and the
# replacement text also happens to look numeric ("836.00"), so a direct
# Range.Value assignment would make Excel re-interpret it as a number and
# drop the trailing zeros. Instead, build the literal text via a
# formula that evaluates to a string, then paste-special "values only"
# into the target cell - this writes a genuine text cell without
# perturbing its number format/style (unlike the quote-prefix trick).
$xlPasteValues = -4163

$scratch = $ws.Range("Z100")

$scratch.Formula = "=""836.00"""
$scratch.Copy()
$ws.Range("H2").PasteSpecial($xlPasteValues)

$scratch.Formula = "=""8753.00"""
$scratch.Copy()
$ws.Range("H3").PasteSpecial($xlPasteValues)

$scratch.ClearContents()
$excel.CutCopyMode = $false
